# Add files via upload
# Populates the new enforcement-action rows (62-72) that were added to the
# "SecuritiesEnforcement" tracker with Outcome / Cause of Action / Civil-or-
# Criminal / Token / Project Name / Blockchain / Amount / 1933-Act /
# 1934-Act / SEC-Office detail that previously only had Date/Case/Description.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 62 - SEC v. Kik Interactive Inc. (Kin)
# ---------------------------------------------------------------------------
$ws.Range("D62").Value = "Ongoing"
$ws.Range("E62").Value = "Unregistered Offering"
$ws.Range("F62").Value = "Civil"
$ws.Range("H62").Value = "Kik Interactive Inc."
$ws.Range("G62").Value = "Kin"
$ws.Range("I62").Value = "Ethereum"
$ws.Range("J62").Value = 100000000
$ws.Range("K62").Value = 1
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = "Washington, D.C."

# ---------------------------------------------------------------------------
# Row 63 - SEC v. Pacheco (IPro Solutions LLC and IPro Network LLC / PRO)
# ---------------------------------------------------------------------------
$ws.Range("D63").Value = "Ongoing"
$ws.Range("E63").Value = "Unregistered Offering and Fraud"
$ws.Range("F63").Value = "Civil"
$ws.Range("H63").Value = "IPro Solutions LLC and IPro Network LLC"
$ws.Range("G63").Value = "PRO"
$ws.Range("I63").Value = "N/A"
$ws.Range("J63").Value = 26000000
$ws.Range("K63").Value = 1
$ws.Range("L63").Value = 1
$ws.Range("M63").Value = "Los Angeles"

# ---------------------------------------------------------------------------
# Row 64 - NextBlock Global Ltd. and Alex Tapscott
# ---------------------------------------------------------------------------
$ws.Range("D64").Value = "Settlement"
$ws.Range("E64").Value = "Unregistered Offering"
$ws.Range("F64").Value = "Civil"
$ws.Range("H64").Value = "NextBlock Global Ltd."
$ws.Range("G64").Value = "N/A"
$ws.Range("I64").Value = "N/A"
$ws.Range("J64").Value = 18400000
$ws.Range("K64").Value = 1
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = "Washington, D.C."

# ---------------------------------------------------------------------------
# Row 65 - SEC v. Natural Diamonds Investment Co., et al. (Argyle Coin, LLC)
# ---------------------------------------------------------------------------
$ws.Range("D65").Value = "Ongoing"
$ws.Range("E65").Value = "Unregistered Offering"
$ws.Range("F65").Value = "Civil"
$ws.Range("H65").Value = "Argyle Coin, LLC"
$ws.Range("H65").Font.Color = 0
$ws.Range("G65").Value = "Argyle"
$ws.Range("I65").Value = "N/A"
$ws.Range("J65").Value = 30000000
$ws.Range("K65").Value = 1
$ws.Range("L65").Value = 1
$ws.Range("M65").Value = "Miami"

# ---------------------------------------------------------------------------
# Row 66 - Mutual Coin Fund LLC and Usman Majeed
# ---------------------------------------------------------------------------
$ws.Range("D66").Value = "Settlement"
$ws.Range("E66").Value = "Unregistered Offering"
$ws.Range("F66").Value = "Civil"
$ws.Range("H66").Value = "Mutual Coin Fund LLC"
$ws.Range("G66").Value = "N/A"
$ws.Range("I66").Value = "N/A"
$ws.Range("J66").Value = 567000
$ws.Range("K66").Value = 1
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = "Washington, D.C."

# ---------------------------------------------------------------------------
# Row 67 - Gladius Network LLC
# ---------------------------------------------------------------------------
$ws.Range("D67").Value = "Settlement"
$ws.Range("E67").Value = "Unregistered Offering"
$ws.Range("F67").Value = "Civil"
$ws.Range("G67").Value = "GLA"
$ws.Range("H67").Value = "Gladius Network LLC"
$ws.Range("I67").Value = "Ethereum"
$ws.Range("J67").Value = 12700000
$ws.Range("K67").Value = 1
$ws.Range("L67").Value = 1
$ws.Range("M67").Value = "Washington, D.C."

# ---------------------------------------------------------------------------
# Row 68 - CoinAlpha Advisors LLC
# ---------------------------------------------------------------------------
$ws.Range("D68").Value = "Settlement"
$ws.Range("E68").Value = "Unregistered Offering"
$ws.Range("F68").Value = "Civil"
$ws.Range("G68").Value = "GLA"
$ws.Range("H68").Value = "CoinAlpha Advisors LLC"
$ws.Range("I68").Value = "N/A"
$ws.Range("J68").Value = 608491
$ws.Range("K68").Value = 1
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = "N/A"

# ---------------------------------------------------------------------------
# Row 69 - Floyd Mayweather, Jr. (Centra Tech Inc. / CTR)
# ---------------------------------------------------------------------------
$ws.Range("D69").Value = "Settlement"
$ws.Range("E69").Value = "Anti-touting"
$ws.Range("F69").Value = "Civil"
$ws.Range("H69").Value = "Centra Tech Inc."
$ws.Range("G69").Value = "CTR"
$ws.Range("I69").Value = "Ethereum"
$ws.Range("J69").Value = 300000
$ws.Range("K69").Value = 1
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = "New York"

# ---------------------------------------------------------------------------
# Row 70 - Khaled Khaled ("DJ Khaled") (Centra Tech Inc. / CTR)
# ---------------------------------------------------------------------------
$ws.Range("D70").Value = "Settlement"
$ws.Range("E70").Value = "Anti-touting"
$ws.Range("F70").Value = "Civil"
$ws.Range("G70").Value = "CTR"
$ws.Range("H70").Value = "Centra Tech Inc."
$ws.Range("I70").Value = "Ethereum"
$ws.Range("J70").Value = 50000
$ws.Range("K70").Value = 1
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = "Los Angeles"

# ---------------------------------------------------------------------------
# Row 71 - Paragon Coin, Inc. (ParagonCoin)
# ---------------------------------------------------------------------------
$ws.Range("D71").Value = "Settlement"
$ws.Range("E71").Value = "Unregistered Offering"
$ws.Range("F71").Value = "Civil"
$ws.Range("H71").Value = "Paragon Coin Inc."
$ws.Range("I71").Value = "Ethereum"
$ws.Range("J71").Value = 12000000
$ws.Range("K71").Value = 1
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = "Boston"

# ---------------------------------------------------------------------------
# Row 72 - CarrierEQ, Inc., d/b/a Airfox (AirTokens)
# ---------------------------------------------------------------------------
$ws.Range("D72").Value = "Settlement"
$ws.Range("E72").Value = "Unregistered Offering"
$ws.Range("F72").Value = "Civil"
$ws.Range("H72").Value = "Arifox"
$ws.Range("G72").Value = "AirTokens"
$ws.Range("I72").Value = "Ethereum"
$ws.Range("J72").Value = 15000000
$ws.Range("K72").Value = 1
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = "Boston"

# Token symbols filled in last for the Paragon / Airfox pair (matches the
# original author's entry order).
$ws.Range("G71").Value = "ParagonCoin"

# ---------------------------------------------------------------------------
# Column G ("Token") was widened slightly to fit the new values.
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 12.5

# ---------------------------------------------------------------------------
# Leave the selection where the author last left it.
# ---------------------------------------------------------------------------
$ws.Range("D73").Select()
